$wb = $excel.ActiveWorkbook

# Sheet: "Intermediate for Mapping"
$wsMap = $wb.Worksheets.Item("Intermediate for Mapping")
$wsMap.Range("V2").Value = "No Detect Data"
$wsMap.Range("W2").Value = "No Detect Data"
$wsMap.Range("V12").Value = "No Detect Data"
$wsMap.Range("W12").Value = "No Detect Data"

# Sheet: "Intermediate Exhibit"
$wsExh = $wb.Worksheets.Item("Intermediate Exhibit")
$wsExh.Range("G4").Value = "NA"
$wsExh.Range("H4").Value = "NA"
$wsExh.Range("G16").Value = "NA"
$wsExh.Range("H16").Value = "NA"
